$d = $word.ActiveDocument

# --- Alteração 1 (context): collapse the 4 empty centered paragraphs
# (originally paragraphs 2-5) into a single, fully empty paragraph. ---
$r = $d.Range($d.Paragraphs.Item(3).Range.Start, $d.Paragraphs.Item(5).Range.End)
$r.Delete()
$p2 = $d.Paragraphs.Item(2)
$p2.Format.Alignment = 0

# --- Alteração 2: Adicionar nome do trabalho ---
# Append a new bold, centered, 14pt paragraph ("Trabalho Final DAS")
# right after the "CTESP Desenvolvimento de Software" title paragraph.
$last = $d.Paragraphs.Item($d.Paragraphs.Count)
$last.Range.InsertParagraphAfter()
$newRange = $d.Paragraphs.Item($d.Paragraphs.Count).Range
$newRange.Text = "Trabalho Final DAS"
$newRange.Font.Size = 14
$newRange.Font.SizeBi = 14
